# New weekly price record: insert a row at the top of the data (row 10,
# right after the most recent entry already present) and push every
# existing record down by one row. This is the "Fruta / hortaliza,
# semanal" update - a fresh Poroto granado quote for Terminal La Palmera
# de La Serena was added, so the whole history shifts and the sheet grows
# from A1:R115 to A1:R116.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10..115 down to 11..116, leaving a blank row 10 (inherits the
# date-formatted style of the row above it, just like Excel's own
# Insert Row command).
$ws.Rows("10:10").Insert()

# Populate the newly inserted row with the new weekly quote.
$ws.Cells.Item(10, 1).Value  = 8
$ws.Cells.Item(10, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 44959
$ws.Cells.Item(10, 5).Value  = 4
$ws.Cells.Item(10, 6).Value  = 100112030
$ws.Cells.Item(10, 7).Value  = "Poroto granado"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 31500
$ws.Cells.Item(10, 12).Value = 32000
$ws.Cells.Item(10, 13).Value = 31750
$ws.Cells.Item(10, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 16).Value = 1270
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
